# Apply crypto price/volume updates scraped on 2023-04-16.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.642.02"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "2.120.60"
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5256"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4559"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.44"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09124"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("E11").Value = "  +2.05%  "
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "2.120.21"
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("E14").Value = "  +2.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.157"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001176"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "97.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.88%  "
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06690"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.010"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.327"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("D23").Value = "30.707.89"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.358"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("D26").Value = "2.367.61"
$ws.Range("E26").Value = "  +1.38%  "
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.566"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.212"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.664"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.378"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.942"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.882"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02642"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06874"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2333"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.29%  "
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6919"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.260"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6506"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.319"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000368"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +21.47%  "
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07306"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.68%  "
